$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 283 (which holds "ISU" / Sulaymaniyah, Iraq).
# This shifts the existing ISU row to 284 and the SFO row to 285.
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new "NVT" / Timbo, Brazil data.
$ws.Cells.Item(283, 1).Value = "NVT"
$ws.Cells.Item(283, 2).Value = "Timbo, Brazil"
$ws.Cells.Item(283, 3).Value = -26.8251
$ws.Cells.Item(283, 4).Value = -49.2695
$ws.Cells.Item(283, 5).Value = "BR"
$ws.Cells.Item(283, 6).Value = "South America"
$ws.Cells.Item(283, 7).Value = "Timbo"

# Match the style of column A used throughout the "colo" column (bold, bordered, centered).
$ws.Cells.Item(282, 1).Copy()
$ws.Cells.Item(283, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
